# Auto-generated: re-apply the author's re-shuffled 'Top N' tie-break ordering
# across the 2011-2017 KPI sheets (bikesharing_output/kpis.xlsx).
$wb = $excel.ActiveWorkbook

# --- Sheet "2011" (16 cell updates) ---
$ws = $wb.Worksheets.Item("2011")
$ws.Range("E2").Value = "B00514 (21)"
$ws.Range("G2").Value = "B00454 (21)"
$ws.Range("I2").Value = "B00159 (20)"
$ws.Range("K2").Value = "B00557 (20)"
$ws.Range("I3").Value = "B00021 (109)"
$ws.Range("J3").Value = "B00206 (109)"
$ws.Range("K3").Value = "B00120 (108)"
$ws.Range("L3").Value = "B00144 (108)"
$ws.Range("G5").Value = "B00217 (102)"
$ws.Range("H5").Value = "B00078 (102)"
$ws.Range("M5").Value = "B00109 (94)"
$ws.Range("N5").Value = "B00127 (94)"
$ws.Range("I6").Value = "B00311 (73)"
$ws.Range("J6").Value = "B00121 (73)"
$ws.Range("K6").Value = "B00353 (72)"
$ws.Range("M6").Value = "B00431 (72)"

# --- Sheet "2012" (13 cell updates) ---
$ws = $wb.Worksheets.Item("2012")
$ws.Range("G2").Value = "B00149 (60)"
$ws.Range("H2").Value = "B00125 (60)"
$ws.Range("I2").Value = "B00008 (59)"
$ws.Range("J2").Value = "B00389 (59)"
$ws.Range("K2").Value = "B00101 (59)"
$ws.Range("M2").Value = "B00221 (57)"
$ws.Range("J5").Value = "B00348 (137)"
$ws.Range("K5").Value = "B00397 (137)"
$ws.Range("J7").Value = "T01195 (173)"
$ws.Range("K7").Value = "T01013 (173)"
$ws.Range("E10").Value = "T01018 (102)"
$ws.Range("F10").Value = "T01158 (102)"
$ws.Range("N10").Value = "T01288 (86)"

# --- Sheet "2013" (12 cell updates) ---
$ws = $wb.Worksheets.Item("2013")
$ws.Range("L3").Value = "T01362 (151)"
$ws.Range("M3").Value = "T01177 (151)"
$ws.Range("L5").Value = "T01068 (183)"
$ws.Range("M5").Value = "B00083 (183)"
$ws.Range("N6").Value = "T01238 (213)"
$ws.Range("E7").Value = "T01283 (201)"
$ws.Range("F7").Value = "T01297 (201)"
$ws.Range("G7").Value = "B01479 (201)"
$ws.Range("N7").Value = "T01225 (187)"
$ws.Range("L10").Value = "T01232 (53)"
$ws.Range("M10").Value = "T01146 (53)"
$ws.Range("N10").Value = "T01393 (53)"

# --- Sheet "2014" (12 cell updates) ---
$ws = $wb.Worksheets.Item("2014")
$ws.Range("I2").Value = "T01255 (46)"
$ws.Range("J2").Value = "T01267 (46)"
$ws.Range("K2").Value = "B01484 (46)"
$ws.Range("L3").Value = "T01013 (46)"
$ws.Range("M3").Value = "B01465 (46)"
$ws.Range("I5").Value = "T01042 (130)"
$ws.Range("J5").Value = "T01152 (130)"
$ws.Range("F9").Value = "B00585 (212)"
$ws.Range("G9").Value = "B01660 (212)"
$ws.Range("N10").Value = "T01093 (194)"
$ws.Range("G11").Value = "B01795 (160)"
$ws.Range("H11").Value = "A07813 (160)"

# --- Sheet "2015" (17 cell updates) ---
$ws = $wb.Worksheets.Item("2015")
$ws.Range("V3").Value = "822 (29)"
$ws.Range("W3").Value = "1139 (29)"
$ws.Range("Q4").Value = "1302 (61)"
$ws.Range("R4").Value = "660 (61)"
$ws.Range("W4").Value = "769 (58)"
$ws.Range("T5").Value = "782 (109)"
$ws.Range("U5").Value = "709 (109)"
$ws.Range("V5").Value = "1289 (108)"
$ws.Range("W5").Value = "948 (108)"
$ws.Range("R6").Value = "774 (182)"
$ws.Range("T6").Value = "614 (182)"
$ws.Range("R7").Value = "1230 (180)"
$ws.Range("S7").Value = "617 (180)"
$ws.Range("W7").Value = "351 (171)"
$ws.Range("X7").Value = "462 (171)"
$ws.Range("R13").Value = "1418 (152)"
$ws.Range("S13").Value = "1495 (152)"

# --- Sheet "2016" (20 cell updates) ---
$ws = $wb.Worksheets.Item("2016")
$ws.Range("Q2").Value = "1565 (63)"
$ws.Range("R2").Value = "654 (63)"
$ws.Range("S2").Value = "643 (59)"
$ws.Range("T2").Value = "681 (59)"
$ws.Range("V2").Value = "1550 (57)"
$ws.Range("W2").Value = "926 (57)"
$ws.Range("O3").Value = "1545 (74)"
$ws.Range("P3").Value = "653 (74)"
$ws.Range("U6").Value = "1443 (224)"
$ws.Range("V6").Value = "1384 (224)"
$ws.Range("X6").Value = "1428 (216)"
$ws.Range("S8").Value = "1385 (241)"
$ws.Range("T8").Value = "1479 (241)"
$ws.Range("X8").Value = "1433 (235)"
$ws.Range("AC8").Value = "Charles Circle - Charles St. at Cambridge St. (2625)"
$ws.Range("AD8").Value = "Boston Public Library - 700 Boylston St. (2625)"
$ws.Range("T11").Value = "1810 (195)"
$ws.Range("U11").Value = "1700 (195)"
$ws.Range("U12").Value = "1880 (155)"
$ws.Range("V12").Value = "1429 (155)"

# --- Sheet "2017" (17 cell updates) ---
$ws = $wb.Worksheets.Item("2017")
$ws.Range("S2").Value = "1559 (68)"
$ws.Range("T2").Value = "1656 (68)"
$ws.Range("U4").Value = "769 (74)"
$ws.Range("V4").Value = "1483 (74)"
$ws.Range("R5").Value = "1434 (148)"
$ws.Range("S5").Value = "277 (148)"
$ws.Range("U6").Value = "1711 (165)"
$ws.Range("V6").Value = "1710 (165)"
$ws.Range("X6").Value = "1811 (165)"
$ws.Range("Q7").Value = "1865 (211)"
$ws.Range("R7").Value = "1828 (211)"
$ws.Range("S10").Value = "1824 (209)"
$ws.Range("T10").Value = "1257 (209)"
$ws.Range("T12").Value = "1900 (160)"
$ws.Range("U12").Value = "1749 (160)"
$ws.Range("U13").Value = "1906 (88)"
$ws.Range("V13").Value = "1893 (88)"

